# Adding new link to scraping
# Appends 6 new rows (234-239) of daily COVID-19 statistics data for
# Bosnia and Herzegovina (bih) to the bottom of the existing data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, Date (col A, kept as literal text, not an Excel date),
# Confirmed cases (B), Tested (C), Deaths (D), Recovered (E), Active (F), Monitored (G)
$newRows = @(
    @(234, "11.01.2021", 115758, 544432, 4358, 82191, 29209, 0),
    @(235, "10.01.2021", 115633, 541875, 4330, 81624, 29679, 0),
    @(236, "09.01.2021", 115379, 538848, 4305, 81325, 29749, 0),
    @(237, "08.01.2021", 114920, 535439, 4285, 80868, 29767, 0),
    @(238, "05.01.2021", 113392, 524907, 4211, 79465, 29716, 0),
    @(239, "03.01.2021", 112645, 519854, 4131, 77891, 30623, 0)
)

foreach ($entry in $newRows) {
    $r = $entry[0]

    # Column A holds a date-like string (e.g. "11.01.2021"). Excel would
    # normally auto-convert such text into a real date serial number, so
    # force the cell to Text format first, assign the literal string, then
    # clear the formatting again so the cell ends up with the default style
    # (no explicit style index) while keeping the value as plain text.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $entry[1]
    $cellA.ClearFormats()

    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]
    $ws.Cells.Item($r, 4).Value = $entry[4]
    $ws.Cells.Item($r, 5).Value = $entry[5]
    $ws.Cells.Item($r, 6).Value = $entry[6]
    $ws.Cells.Item($r, 7).Value = $entry[7]
}
